$wb = $excel.ActiveWorkbook

# --- NoPowerState: quicker discharge cap check threshold raised to 10V, update selection ---
$wsNoPower = $wb.Worksheets.Item("NoPowerState")
$wsNoPower.Activate() | Out-Null
$wsNoPower.Range("C14").Value = 10
$wsNoPower.Range("H7").Select() | Out-Null

# --- NormalOperationState: TEMP limits row (B7:D7) lowered since board was heated too long ---
$wsNormal = $wb.Worksheets.Item("NormalOperationState")
$wsNormal.Activate() | Out-Null
$wsNormal.Range("B7").Value = 23
$wsNormal.Range("C7").Value = 28
$wsNormal.Range("D7").Value = 25
$wsNormal.Range("C7").Select() | Out-Null

# --- CapsChargingState: TP2B row (B4:D4) updated temps ---
$wsCaps = $wb.Worksheets.Item("CapsChargingState")
$wsCaps.Activate() | Out-Null
$wsCaps.Range("B4").Value = 23.8
$wsCaps.Range("C4").Value = 26.2
$wsCaps.Range("D4").Value = 25
$wsCaps.Range("F4").Select() | Out-Null

# --- SPMState: TEMP limits row (B7:D7) lowered, matches NormalOperationState change ---
$wsSPM = $wb.Worksheets.Item("SPMState")
$wsSPM.Activate() | Out-Null
$wsSPM.Range("B7").Value = 23
$wsSPM.Range("C7").Value = 28
$wsSPM.Range("D7").Value = 25
$wsSPM.Range("B7").Select() | Out-Null

# --- Report: keep cursor at D29, no longer the active tab ---
$wsReport = $wb.Worksheets.Item("Report")
$wsReport.Activate() | Out-Null
$wsReport.Range("D29").Select() | Out-Null

# --- Quantities becomes the active/visible sheet ---
$wsQuantities = $wb.Worksheets.Item("Quantities")
$wsQuantities.Activate() | Out-Null
$wsQuantities.Range("I11").Select() | Out-Null
